$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5680
$ws.Range("B1").Value = 5680
$ws.Range("C1").Value = 5773.999999999985

$ws.Range("A2").Value = 5680
$ws.Range("B2").Value = 5680
$ws.Range("C2").Value = 5778

$ws.Range("A3").Value = 5680
$ws.Range("B3").Value = 5680
$ws.Range("C3").Value = 5733.999999999997

$ws.Range("A4").Value = 5680
$ws.Range("B4").Value = 11552.99411623487
$ws.Range("C4").Value = 5680

$ws.Range("A5").Value = 5680
$ws.Range("B5").Value = 5767.999999999498
$ws.Range("C5").Value = 5680

$ws.Range("A6").Value = 5680
$ws.Range("B6").Value = 11552.99411623487
$ws.Range("C6").Value = 5680
